$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Auto-update: append the newest Pick 3 draw result as row 69.
# A (date) and C (6-digit phase code) look numeric/date-like to Excel's
# smart-input parser, so they are entered with a leading apostrophe to keep
# them as literal text (matches the rest of the column's text values).
$ws.Range("A69").Value = "'2025-11-24"
$ws.Range("B69").Value = "Pick 3"
$ws.Range("C69").Value = "'251124"
$ws.Range("D69").Value = "2-4-2"
$ws.Range("E69").Value = "2025-11-24T21:41:01.577+04:00"
